$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Bmp7", "Bmpr1a", "ECs", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 6.322177333333333, 18.966532, 0.08271011762055308, 0.08271011762055309, 0.2177315725751111, 1.959584153176, 0.004816426702708783, 0.004816426702708785),
    @("ECs", "Bmp7", "Bmpr1a", "FAPs", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 41.286995, 123.860985, 0.5401386314560596, 0.5401386314560597, 1.421896583136666, 12.79706924823, 0.03145368671393442, 0.03145368671393443),
    @("ECs", "Bmp7", "Bmpr1a", "MuSCs", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 27.73243066666667, 83.197292, 0.3628105447549136, 0.3628105447549136, 0.9550864238728889, 8.595777814856, 0.0211274079405692, 0.0211274079405692),
    @("ECs", "Bmp7", "Bmpr1a", "Resolving-Mac", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 1.096171666666667, 3.288515, 0.01434070616847367, 0.01434070616847367, 0.03775142141888889, 0.33976279277, 0.0008350968673797812, 0.0008350968673797813),
    @("FAPs", "Bmp7", "Bmpr1a", "ECs", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 6.322177333333333, 18.966532, 0.08271011762055308, 0.08271011762055309, 3.263824048333333, 29.374416435, 0.07219885069223396, 0.07219885069223399),
    @("FAPs", "Bmp7", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 41.286995, 123.860985, 0.5401386314560596, 0.5401386314560597, 21.31441116875, 191.82970051875, 0.4714947763042833, 0.4714947763042834),
    @("FAPs", "Bmp7", "Bmpr1a", "MuSCs", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 27.73243066666667, 83.197292, 0.3628105447549136, 0.3628105447549136, 14.31686733166667, 128.851805985, 0.3167025401958667, 0.3167025401958667),
    @("FAPs", "Bmp7", "Bmpr1a", "Resolving-Mac", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 1.096171666666667, 3.288515, 0.01434070616847367, 0.01434070616847367, 0.5658986229166667, 5.093087606250001, 0.01251820857308926, 0.01251820857308926),
    @("MuSCs", "Bmp7", "Bmpr1a", "ECs", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 6.322177333333333, 18.966532, 0.08271011762055308, 0.08271011762055309, 0.2574411684057778, 2.316970515652, 0.005694840225610327, 0.005694840225610328),
    @("MuSCs", "Bmp7", "Bmpr1a", "FAPs", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 41.286995, 123.860985, 0.5401386314560596, 0.5401386314560597, 1.681220198731667, 15.130981788585, 0.03719016843784183, 0.03719016843784185),
    @("MuSCs", "Bmp7", "Bmpr1a", "MuSCs", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 27.73243066666667, 83.197292, 0.3628105447549136, 0.3628105447549136, 1.129273820890222, 10.163464388012, 0.02498059661847766, 0.02498059661847766),
    @("MuSCs", "Bmp7", "Bmpr1a", "Resolving-Mac", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 1.096171666666667, 3.288515, 0.01434070616847367, 0.01434070616847367, 0.04463647565722223, 0.4017282809150001, 0.0009874007280046213, 0.0009874007280046213),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
